$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous event (row 11) failed to upload correctly. Re-create it in
# row 12 with the same captured data, fixing the ID so the new upload is
# recognised as a distinct event.
$ws.Range("A11:U11").Copy($ws.Range("A12:U12"))
$ws.Range("A12").Value = "5c5479fbe394015b5e3941e7"

# A few of the copied fields are blank (empty text) answers. A plain
# assignment of "" clears the cell entirely instead of leaving an empty
# text value in it, so use the quote-prefix trick to force an explicit
# (but empty) text entry, then drop the resulting quote-prefix formatting
# so the cell keeps the workbook's default style.
foreach ($col in @("I", "O", "P", "S", "T", "U")) {
    $cell = $ws.Range($col + "12")
    $cell.Value = "'"
    $cell.Style = "Normal"
}

# Move the active selection to the newly filled-in row, like a user would
# after finishing the edit.
$ws.Range("A12").Select()
